# Edit: add "version" and "description" columns to "Export as TSV" sheet,
# and add a new "version list" sheet (inserted before "assay_category list").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert two new columns (A:B) on the main sheet. This shifts all
#        existing data, formulas, and data-validation ranges right by two
#        columns automatically. ---
$ws.Columns("A:B").Insert()

# --- 2. Apply the same header formatting as the rest of row 1 (bold,
#        centered, wrapped) to the two new header cells, then set their
#        text. ---
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("A1").Value = "version"
$ws.Range("B1").Value = "description"

# --- 3. Re-create the header comments. Inserting columns does not move
#        existing cell comments, so remove them all and re-add them at
#        their new (shifted) locations, plus the two new ones for the
#        "version" and "description" columns. ---
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}

$ws.Range("A1").AddComment("Version of the schema to use when validating this metadata.") | Out-Null
$ws.Range("B1").AddComment("Free-text description of this assay.") | Out-Null

$shiftedComments = @(
    @{Ref='C1'; Text='HuBMAP Display ID of the donor of the assayed tissue.'}
    @{Ref='D1'; Text='HuBMAP Display ID of the assayed tissue.'}
    @{Ref='E1'; Text='Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.'}
    @{Ref='F1'; Text='DOI for protocols.io referring to the protocol for this assay.'}
    @{Ref='G1'; Text='Name of the person responsible for executing the assay.'}
    @{Ref='H1'; Text='Email address for the operator.'}
    @{Ref='I1'; Text='Name of the principal investigator responsible for the data.'}
    @{Ref='J1'; Text='Email address for the principal investigator.'}
    @{Ref='K1'; Text='Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.'}
    @{Ref='L1'; Text='The specific type of assay being executed.'}
    @{Ref='M1'; Text='Analytes are the target molecules being measured with the assay.'}
    @{Ref='N1'; Text='Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.'}
    @{Ref='O1'; Text='An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.'}
    @{Ref='P1'; Text='Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.'}
    @{Ref='Q1'; Text='A number (no comma separators)'}
    @{Ref='R1'; Text='Link to a protocols document answering the question: How was tissue stored and processed for cell/nuclei isolation'}
    @{Ref='S1'; Text='Is this a sequencing replicate?'}
    @{Ref='T1'; Text='Adapter sequence to be used for adapter trimming'}
    @{Ref='U1'; Text='Average size in basepairs (bp) of sequencing library fragments estimated via gel electrophoresis or Bioanalyzer/tapestation.'}
    @{Ref='V1'; Text='The concentration value of the pooled library samples submitted for sequencing.'}
    @{Ref='W1'; Text='Unit of library_concentration_value'}
    @{Ref='X1'; Text='A link to the protocol document containing the library construction method (including version) that was used, e.g. "Smart-Seq2", "Drop-Seq", "10X v3".'}
    @{Ref='Y1'; Text='date and time of library creation. YYYY-MM-DD, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s.'}
    @{Ref='Z1'; Text='Total amount (eg. nanograms) of library after the clean-up step of final pcr amplification step. Answer the question: What is the Qubit measured concentration (ng/ul) times the elution volume (ul) after the final clean-up step?'}
    @{Ref='AA1'; Text='Units of library final yield'}
    @{Ref='AB1'; Text='ID of the library sample.'}
    @{Ref='AC1'; Text='State whether the library was generated for single-end or paired end sequencing.'}
    @{Ref='AD1'; Text='Number of PCR cycles performed in order to add adapters and amplify the library. Usually, this includes 5 pre-amplificationn cycles followed by 0-5 additional cycles determined by qPCR.'}
    @{Ref='AE1'; Text='Reagent kit used for library preparation'}
    @{Ref='AF1'; Text='This is a quality metric by visual inspection. This should answer the question: Are the nuclei intact and are the nuclei free of significant amounts of debris? This can be captured at a high level, “OK” or “not OK”.'}
    @{Ref='AG1'; Text='Percent PhiX loaded to the run'}
    @{Ref='AH1'; Text='Slash-delimited list of the number of sequencing cycles for, for example, Read1, i7 index, i5 index, and Read2.'}
    @{Ref='AI1'; Text='Percent of bases with Quality scores above Q30'}
    @{Ref='AJ1'; Text='Reagent kit used for sequencing. NovaSeq6000 for example'}
    @{Ref='AK1'; Text='If Tn5 came from a kit, provide the catalog number.'}
    @{Ref='AL1'; Text='Modality of capturing accessible chromatin molecules. The kit used, for example.'}
    @{Ref='AM1'; Text='The source of the Tn5 transposase and transposon used for capturing accessible chromatin.'}
    @{Ref='AN1'; Text='Relative path to file with ORCID IDs for contributors for this dataset.'}
    @{Ref='AO1'; Text='Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.'}
)

foreach ($item in $shiftedComments) {
    $ws.Range($item.Ref).AddComment($item.Text) | Out-Null
}

# --- 4. Insert the new "version list" sheet before "assay_category list"
#        (which is currently the second sheet), and give it the single
#        allowed value for the new "version" column. ---
$listSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$listSheet.Name = "version list"
$listSheet.Range("A1").NumberFormat = "@"
$listSheet.Range("A1").Value = "1"
$listSheet.Range("A1").Style = "Normal"

# --- 5. Add data validation to the new "version" column, restricting it
#        to the values from the "version list" sheet (mirroring the
#        pattern used for the other list-backed columns). ---
$rng = $ws.Range("A2:A1048576")
$rng.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1") | Out-Null
$rng.Validation.ErrorTitle = "Value must come from list"
$rng.Validation.ErrorMessage = "Value must be one of: 1."
$rng.Validation.IgnoreBlank = $true
$rng.Validation.InCellDropdown = $true
